# Update unit tests for existing catalog
# - Rename Sheet1 -> Catalog, add a new Properties sheet
# - Catalog: insert "Base Directory"/"Relative Path" columns, rename
#   "File Size (bytes)" -> "Readable Size", add "Checksum"/"File Link"/
#   "Link Directory"/"Link Name" columns, and refresh the sample rows
# - Properties: new sheet describing the catalog run configuration

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Catalog"

# ---- Catalog header row (B1:N1) ----
$ws.Range("B1").Value = "File Path"
$ws.Range("C1").Value = "Base Directory"
$ws.Range("D1").Value = "Relative Path"
$ws.Range("E1").Value = "Filename"
$ws.Range("F1").Value = "Extension"
$ws.Range("G1").Value = "File Size"
$ws.Range("H1").Value = "Readable Size"
$ws.Range("I1").Value = "Checksum"
$ws.Range("J1").Value = "Duplicate"
$ws.Range("K1").Value = "File Link"
$ws.Range("L1").Value = "Directory"
$ws.Range("M1").Value = "Link Directory"
$ws.Range("N1").Value = "Link Name"

# New header cells J1:N1 need the same bold/centered/bordered look as B1:I1
$ws.Range("B1").Copy()
$ws.Range("J1:N1").PasteSpecial(-4122)

# ---- Catalog data rows (A2:J7) ----
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "C:\Users\cdevore\Documents\GitHub\DocumentCatalog\test\Email Security Digest 5 New Messages .msg"
$ws.Range("C2").Value = "C:\Users\cdevore\Documents\GitHub\DocumentCatalog\test"
$ws.Range("D2").Value = "Email Security Digest 5 New Messages .msg"
$ws.Range("E2").Value = "Email Security Digest 5 New Messages .msg"
$ws.Range("F2").Value = ".msg"
$ws.Range("G2").Value = 65024
$ws.Range("H2").Value = "64KB"
$ws.Range("I2").Value = "b8b8f59c500d3ce9e6392e1c1b2ffc53af78e838"
$ws.Range("J2").Value = $false

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "C:\Users\cdevore\Documents\GitHub\DocumentCatalog\test\email02.msg"
$ws.Range("C3").Value = "C:\Users\cdevore\Documents\GitHub\DocumentCatalog\test"
$ws.Range("D3").Value = "email02.msg"
$ws.Range("E3").Value = "email02.msg"
$ws.Range("F3").Value = ".msg"
$ws.Range("G3").Value = 66048
$ws.Range("H3").Value = "64KB"
$ws.Range("I3").Value = "daa063c933cbdfd82dae57b451dcc488c8c19a0f"
$ws.Range("J3").Value = $false

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "C:\Users\cdevore\Documents\GitHub\DocumentCatalog\test\some_files.xlsx"
$ws.Range("C4").Value = "C:\Users\cdevore\Documents\GitHub\DocumentCatalog\test"
$ws.Range("D4").Value = "some_files.xlsx"
$ws.Range("E4").Value = "some_files.xlsx"
$ws.Range("F4").Value = ".xlsx"
$ws.Range("G4").Value = 6085
$ws.Range("H4").Value = "6KB"
$ws.Range("I4").Value = "187654a5831d5fec4c497a59f78d4c13aae7fffc"
$ws.Range("J4").Value = $false

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "C:\Users\cdevore\Documents\GitHub\DocumentCatalog\test\text1.txt"
$ws.Range("C5").Value = "C:\Users\cdevore\Documents\GitHub\DocumentCatalog\test"
$ws.Range("D5").Value = "text1.txt"
$ws.Range("E5").Value = "text1.txt"
$ws.Range("F5").Value = ".txt"
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = "0B"
$ws.Range("I5").Value = "da39a3ee5e6b4b0d3255bfef95601890afd80709"
$ws.Range("J5").Value = $false

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "C:\Users\cdevore\Documents\GitHub\DocumentCatalog\test\this is an Excel spreadsheet.xlsx"
$ws.Range("C6").Value = "C:\Users\cdevore\Documents\GitHub\DocumentCatalog\test"
$ws.Range("D6").Value = "this is an Excel spreadsheet.xlsx"
$ws.Range("E6").Value = "this is an Excel spreadsheet.xlsx"
$ws.Range("F6").Value = ".xlsx"
$ws.Range("G6").Value = 6171
$ws.Range("H6").Value = "6KB"
$ws.Range("I6").Value = "edf0ec1ae0430ed567294e292a69dd371a4de939"
$ws.Range("J6").Value = $false

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "C:\Users\cdevore\Documents\GitHub\DocumentCatalog\test\this_is_a_word_document.docx"
$ws.Range("C7").Value = "C:\Users\cdevore\Documents\GitHub\DocumentCatalog\test"
$ws.Range("D7").Value = "this_is_a_word_document.docx"
$ws.Range("E7").Value = "this_is_a_word_document.docx"
$ws.Range("F7").Value = ".docx"
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = "0B"
$ws.Range("I7").Value = "da39a3ee5e6b4b0d3255bfef95601890afd80709"
$ws.Range("J7").Value = $true

# ---- New "Properties" sheet, placed after Catalog ----
$props = $wb.Worksheets.Add($null, $ws)
$props.Name = "Properties"

$props.Range("A1").Value = "Document Catalog Properties"
$props.Range("A3").Value = "Search Directories:"
$props.Range("B3").Value = "C:\Users\cdevore\Documents\GitHub\DocumentCatalog\test"
$props.Range("A4").Value = "Exclude Directories:"
$props.Range("B4").Value = "sub_dir"
$props.Range("A5").Value = "Buffer Size:"
$props.Range("B5").Value = 65536
$props.Range("A6").Value = "Hash Function:"
$props.Range("B6").Value = "sha1"

# Restore the Catalog tab as the active/selected sheet
$ws.Activate()
